$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "v1" row values
$ws.Range("E2").Value = 0.35
$ws.Range("F2").Value = 4.46
$ws.Range("G2").Value = 123.94
$ws.Range("H2").Value = 128.76

# Add new "v2" row
$ws.Range("A3").Value = "v2"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.1
$ws.Range("F3").Value = 2.37
$ws.Range("G3").Value = 55.12
$ws.Range("H3").Value = 57.6

# Match Excel's end-of-edit active cell/selection state
$ws.Range("H3").Select() | Out-Null
